# Auto-generated script applying scheduled market-price refresh to Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 54.875
$ws.Range("I38").Value = 54.875
$ws.Range("K38").Value = 164.625
$ws.Range("M38").Value = 207.375
$ws.Range("H58").Value = 343.33334
$ws.Range("I58").Value = 343.33334
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1030.00002
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -880.0000199999999
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 800
$ws.Range("I61").Value = 800
$ws.Range("K61").Value = 2400
$ws.Range("M61").Value = -2228
$ws.Range("H62").Value = 9583
$ws.Range("I62").Value = 7962.6665
$ws.Range("K62").Value = 7962.6665
$ws.Range("M62").Value = -7338.6665
$ws.Range("H64").Value = 9999.454
$ws.Range("I64").Value = 6499.375
$ws.Range("J64").Value = 19333
$ws.Range("K64").Value = 6499.375
$ws.Range("L64").Value = 19333
$ws.Range("M64").Value = -6251.375
$ws.Range("N64").Value = -19829
$ws.Range("H65").Value = 9583
$ws.Range("I65").Value = 7962.6665
$ws.Range("K65").Value = 39813.3325
$ws.Range("M65").Value = -36693.3325
$ws.Range("H67").Value = 9999.454
$ws.Range("I67").Value = 6499.375
$ws.Range("J67").Value = 19333
$ws.Range("K67").Value = 6499.375
$ws.Range("L67").Value = 19333
$ws.Range("M67").Value = -5641.375
$ws.Range("N67").Value = -21049
$ws.Range("H76").Value = 2335
$ws.Range("I76").Value = 2468.75
$ws.Range("J76").Value = 1800
$ws.Range("K76").Value = 2468.75
$ws.Range("L76").Value = 1800
$ws.Range("M76").Value = -2153.75
$ws.Range("N76").Value = -2430
$ws.Range("H79").Value = 2335
$ws.Range("I79").Value = 2468.75
$ws.Range("J79").Value = 1800
$ws.Range("K79").Value = 2468.75
$ws.Range("L79").Value = 1800
$ws.Range("M79").Value = -1376.75
$ws.Range("N79").Value = -3984
$ws.Range("H100").Value = 2960.75
$ws.Range("I100").Value = 3181.7646
$ws.Range("J100").Value = 1708.3334
$ws.Range("K100").Value = 3181.7646
$ws.Range("L100").Value = 1708.3334
$ws.Range("M100").Value = -2640.7646
$ws.Range("N100").Value = -2790.3334
$ws.Range("H103").Value = 4459.3335
$ws.Range("I103").Value = 3000
$ws.Range("J103").Value = 6405.1113
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 19215.3339
$ws.Range("M103").Value = -8414
$ws.Range("N103").Value = -20387.3339
$ws.Range("H113").Value = 14250
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 25000
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 25000
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -31508
$ws.Range("H115").Value = 1104.5714
$ws.Range("I115").Value = 1104.5714
$ws.Range("K115").Value = 3313.7142
$ws.Range("M115").Value = -1746.7142
$ws.Range("H132").Value = 4183.75
$ws.Range("I132").Value = 2880
$ws.Range("J132").Value = 5487.5
$ws.Range("K132").Value = 8640
$ws.Range("L132").Value = 16462.5
$ws.Range("M132").Value = -6110
$ws.Range("N132").Value = -21522.5
$ws.Range("H137").Value = 3490.9092
$ws.Range("I137").Value = 950
$ws.Range("J137").Value = 4055.5557
$ws.Range("K137").Value = 2850
$ws.Range("L137").Value = 12166.6671
$ws.Range("M137").Value = -300
$ws.Range("N137").Value = -17266.6671
$ws.Range("H138").Value = 8867.868
$ws.Range("I138").Value = 6782.3335
$ws.Range("J138").Value = 9515.103999999999
$ws.Range("K138").Value = 20347.0005
$ws.Range("L138").Value = 28545.312
$ws.Range("M138").Value = -15207.0005
$ws.Range("N138").Value = -38825.312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4820.415
$ws.Range("I32").Value = 3785.3674
$ws.Range("K32").Value = 3785.3674
$ws.Range("M32").Value = -3498.3674
$ws.Range("H74").Value = 2849
$ws.Range("I74").Value = 2849
$ws.Range("K74").Value = 2849
$ws.Range("M74").Value = -1975
$ws.Range("H77").Value = 2849
$ws.Range("I77").Value = 2849
$ws.Range("K77").Value = 14245
$ws.Range("M77").Value = -9877

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4643
$ws.Range("I20").Value = 3607.3333
$ws.Range("J20").Value = 7750
$ws.Range("K20").Value = 3607.3333
$ws.Range("L20").Value = 7750
$ws.Range("M20").Value = -3360.3333
$ws.Range("N20").Value = -8244
$ws.Range("H94").Value = 1431.3636
$ws.Range("I94").Value = 845.0625
$ws.Range("J94").Value = 2994.8333
$ws.Range("K94").Value = 845.0625
$ws.Range("L94").Value = 2994.8333
$ws.Range("M94").Value = -394.0625
$ws.Range("N94").Value = -3896.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 4009
$ws.Range("I23").Value = 4009
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 4009
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3769
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 4009
$ws.Range("I27").Value = 4009
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 4009
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3817
$ws.Range("N27").ClearContents()
$ws.Range("H29").Value = 5500
$ws.Range("J29").Value = 5500
$ws.Range("L29").Value = 5500
$ws.Range("N29").Value = -6086
$ws.Range("H31").Value = 3038.9167
$ws.Range("I31").Value = 2863.6
$ws.Range("J31").Value = 3915.5
$ws.Range("K31").Value = 2863.6
$ws.Range("L31").Value = 3915.5
$ws.Range("M31").Value = -2568.6
$ws.Range("N31").Value = -4505.5
$ws.Range("H34").Value = 3038.9167
$ws.Range("I34").Value = 2863.6
$ws.Range("J34").Value = 3915.5
$ws.Range("K34").Value = 2863.6
$ws.Range("L34").Value = 3915.5
$ws.Range("M34").Value = -2661.6
$ws.Range("N34").Value = -4319.5
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H69").Value = 6837.636
$ws.Range("I69").Value = 5514.5
$ws.Range("K69").Value = 5514.5
$ws.Range("M69").Value = -4765.5
$ws.Range("H72").Value = 6837.636
$ws.Range("I72").Value = 5514.5
$ws.Range("K72").Value = 16543.5
$ws.Range("M72").Value = -12799.5
$ws.Range("H114").Value = 105995
$ws.Range("J114").Value = 105995
$ws.Range("L114").Value = 105995
$ws.Range("N114").Value = -114673
$ws.Range("H141").Value = 381258.6
$ws.Range("J141").Value = 381258.6
$ws.Range("L141").Value = 381258.6
$ws.Range("N141").Value = -391618.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 925.8
$ws.Range("I137").Value = 907.5
$ws.Range("J137").Value = 999
$ws.Range("K137").Value = 2722.5
$ws.Range("L137").Value = 2997
$ws.Range("M137").Value = 2377.5
$ws.Range("N137").Value = -13197

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 14166.667
$ws.Range("J29").Value = 14166.667
$ws.Range("L29").Value = 14166.667
$ws.Range("N29").Value = -14746.667
$ws.Range("H43").Value = 2899.5
$ws.Range("I43").Value = 2899.5
$ws.Range("K43").Value = 2899.5
$ws.Range("M43").Value = -2748.5
$ws.Range("H132").Value = 4418
$ws.Range("I132").Value = 4418
$ws.Range("K132").Value = 13254
$ws.Range("M132").Value = -10724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H61").Value = 2247.8462
$ws.Range("I61").Value = 1742.4286
$ws.Range("K61").Value = 1742.4286
$ws.Range("M61").Value = -1540.4286
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H113").Value = 2247.8462
$ws.Range("I113").Value = 1742.4286
$ws.Range("K113").Value = 1742.4286
$ws.Range("M113").Value = 427.5714
$ws.Range("H132").Value = 5630
$ws.Range("I132").Value = 5630
$ws.Range("K132").Value = 16890
$ws.Range("M132").Value = -14360
$ws.Range("H136").Value = 8332.666999999999
$ws.Range("I136").Value = 4998
$ws.Range("K136").Value = 14994
$ws.Range("M136").Value = -12444
$ws.Range("H139").Value = 195000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 539
$ws.Range("I21").Value = 539
$ws.Range("K21").Value = 539
$ws.Range("M21").Value = -304
$ws.Range("H35").Value = 539
$ws.Range("I35").Value = 539
$ws.Range("K35").Value = 539
$ws.Range("M35").Value = -249
$ws.Range("H81").Value = 1259.8
$ws.Range("I81").Value = 1433
$ws.Range("K81").Value = 2866
$ws.Range("M81").Value = -1805
$ws.Range("H84").Value = 1259.8
$ws.Range("I84").Value = 1433
$ws.Range("K84").Value = 14330
$ws.Range("M84").Value = -9026
$ws.Range("H96").Value = 2100
$ws.Range("I96").Value = 1616.6666
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 1616.6666
$ws.Range("L96").Value = 5000
$ws.Range("M96").Value = -243.6666
$ws.Range("N96").Value = -7746
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680
$ws.Range("H122").Value = 753.2
$ws.Range("I122").Value = 754
$ws.Range("K122").Value = 2262
$ws.Range("M122").Value = 188
$ws.Range("H136").Value = 2192.1667
$ws.Range("I136").Value = 1635.3043
$ws.Range("K136").Value = 4905.9129
$ws.Range("M136").Value = -2355.9129

Write-Host "Applied scheduled price update."
